$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10-19 (shrinks used range to A1:H9)
$ws.Range("A10:A19").EntireRow.Delete()

# Adjust column widths: B -> 48, H -> 12 (ColumnWidth = xml width - 0.83)
$ws.Columns.Item(2).ColumnWidth = 47.17
$ws.Columns.Item(8).ColumnWidth = 11.17

# Rewrite rows 2-9 with the freshly scraped listings

# Row 2
$ws.Cells.Item(2, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(2, 2).Value = "AIオートメーションエンジニア"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5452520"
$ws.Cells.Item(2, 7).Value = 303
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai"

# Row 3
$ws.Cells.Item(3, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(3, 2).Value = "【Flutterエンジニア募集】Androidアプリ開発のパートナーを探しています"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5452211"
$ws.Cells.Item(3, 7).Value = 100
$ws.Cells.Item(3, 8).Value = "◆開発 ◇アプリ"

# Row 4
$ws.Cells.Item(4, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(4, 2).Value = "【再掲】基幹システム入替に伴うBIツール環境の再構築(Microsoft Power BI)"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5452367"
$ws.Cells.Item(4, 7).Value = 88
$ws.Cells.Item(4, 8).Value = "◆ツール"

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(5, 2).Value = "Amazonの購入アカウントから必要な情報のスクレイピング→スプレッドシートに記入をしたい。"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5452210"
$ws.Cells.Item(5, 7).Value = 40
$ws.Cells.Item(5, 8).Value = "◆スクレイピング"

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(6, 2).Value = "【R/Shiny】高齢者評価アプリ 機能追加・UI改修依頼"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5452159"
$ws.Cells.Item(6, 7).Value = 38
$ws.Cells.Item(6, 8).Value = "◇アプリ"

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(7, 2).Value = "【小規模・短納期・急募】アプリMatrixifyを用いたデータ移行検証・マッピング担当募集"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5451926"
$ws.Cells.Item(7, 7).Value = 33
$ws.Cells.Item(7, 8).Value = "◇アプリ"

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(8, 2).Value = "注目 限定公開 PR 限定公開の仕事"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5450323"
$ws.Cells.Item(8, 7).Value = 13
$ws.Cells.Item(8, 8).ClearContents()

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-12-12 06:29:31"
$ws.Cells.Item(9, 2).Value = "Xの運用代行"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5451931"
$ws.Cells.Item(9, 7).Value = 10
$ws.Cells.Item(9, 8).ClearContents()

# Rebuild hyperlinks for F2:F9 pointing at the new URLs (old ones were wiped by the row delete/shift)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5452520")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5452211")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5452367")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5452210")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5452159")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5451926")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5450323")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5451931")
